# Update the NATMI TPM numbers for the Pdgfb-Art1 sheet and drop the
# now-redundant duplicate rows (the source data only has 3 sender/receiver
# combinations after the re-run, not 6).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: ECs -> ECs (unchanged pairing, refreshed numbers) ---
$ws.Range("G2").Value = 35.00391733333333
$ws.Range("H2").Value = 105.011752
$ws.Range("I2").Value = 0.9591895364534718
$ws.Range("J2").Value = 0.9591895364534718
$ws.Range("M2").Value = 0.01257466666666667
$ws.Range("N2").Value = 0.037724
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.4401625924942222
$ws.Range("R2").Value = 3.961463332448
$ws.Range("S2").Value = 0.9591895364534718
$ws.Range("T2").Value = 0.9591895364534718

# --- Row 3: was FAPs -> MuSCs, now FAPs -> ECs ---
$ws.Range("A3").Value = "FAPs"
$ws.Range("D3").Value = "ECs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.030543
$ws.Range("H3").Value = 0.091629
$ws.Range("I3").Value = 0.0008369499257158872
$ws.Range("J3").Value = 0.0008369499257158872
$ws.Range("M3").Value = 0.01257466666666667
$ws.Range("N3").Value = 0.037724
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.000384068044
$ws.Range("R3").Value = 0.003456612396
$ws.Range("S3").Value = 0.0008369499257158872
$ws.Range("T3").Value = 0.0008369499257158872

# --- Row 4: was MuSCs -> ECs, now MuSCs -> ECs with refreshed numbers ---
$ws.Range("A4").Value = "MuSCs"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.458762333333333
$ws.Range("H4").Value = 4.376287
$ws.Range("I4").Value = 0.03997351362081222
$ws.Range("J4").Value = 0.03997351362081222
$ws.Range("M4").Value = 0.01257466666666667
$ws.Range("N4").Value = 0.037724
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 0.01834345008755555
$ws.Range("R4").Value = 0.165091050788
$ws.Range("S4").Value = 0.03997351362081222
$ws.Range("T4").Value = 0.03997351362081222

# --- Drop the old rows 5-7 (duplicate combinations no longer present) ---
$ws.Range("A5:T7").Delete()
